# Apply cryptos.xlsx update (Thu Sep 28 12:49:28 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.497.97'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '1.622.84'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.501'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.246'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0609'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.19'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '1.851.73'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '1.612.78'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '235.04'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '26.509.18'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.57%  '
$ws.Range('B24').Value = 'Avalanche'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0495'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').Value = '1.521.39'
$ws.Range('E32').Value = '  +5.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('E35').Value = '  +2.37%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.569'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.835'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  +0.53%  '
$ws.Range('D43').Value = '1.762.78'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.762'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.911'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.88%  '
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  +4.98%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0501'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0965'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.07%  '
